$wb = $excel.ActiveWorkbook

# --- Sheet "组织" (Organization): add a new "组织分类" column (D) ---
$wsOrg = $wb.Worksheets.Item("组织")
$wsOrg.Range("D1").Value = "组织分类"
$wsOrg.Range("D3").Value = "HR组织"
$wsOrg.Range("D4").Value = "HR组织"
$wsOrg.Range("D5").Value = "HR组织"
$wsOrg.Range("D6").Value = "HR组织"
[void]$wsOrg.Range("E5").Select()

# --- Sheet3: just move the selection (no data changes) ---
$wsTodo = $wb.Worksheets.Item("Sheet3")
[void]$wsTodo.Range("G3").Select()

# --- Sheet "人员" (Person): rename header "组织" -> "部门" ---
# Selected/activated last so it remains the workbook's active sheet,
# matching the original file (人员 was the active tab).
$wsPerson = $wb.Worksheets.Item("人员")
$wsPerson.Range("B1").Value = "部门"
[void]$wsPerson.Range("D10").Select()
